$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "37.751.80", "0.634" or "19.50" that
# are stored as literal text in the source data (note the non-numeric double
# "thousands.thousands.decimal" groupings, and values like "19.50"/"4.10" whose
# trailing zero must survive). Temporarily mark the cell as Text before writing
# the value so Excel does not auto-coerce it into a Double, then restore the
# default "Normal" style so no extra formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.751.80'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.084.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.65%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '58.14'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  +2.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.392.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.779'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.71%  '
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.076.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.753.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0836'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E26").Value = '  +8.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.83%  '
$ws.Range("E28").Value = '  -3.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.122'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.97%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0636'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("E37").Value = '  -2.26%  '
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("E40").Value = '  +9.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '101.17'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0966'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("E43").Value = '  +1.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.87%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.450.73'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.68%  '
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.21'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.276.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.41%  '
